$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values as plain text (e.g. "23.211.46", "303.42")
# even when they look numeric, so force a Text number format on each target
# cell before assigning, otherwise Excel would coerce them into real numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.217.35"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "1.601.71"
$ws.Range("E3").Value = "  +0.00%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "1.001"

$ws.Range("D6").Value = "303.70"
$ws.Range("E6").Value = "  +0.71%  "

$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("D8").Value = "52.12"
$ws.Range("E8").Value = "  +4.32%  "

$ws.Range("D9").Value = "0.3631"
$ws.Range("E9").Value = "  -0.26%  "

$ws.Range("D10").Value = "1.271"
$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("D12").Value = "0.08142"
$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("D13").Value = "22.75"
$ws.Range("E13").Value = "  +0.42%  "

$ws.Range("D14").Value = "6.566"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").Value = "7.404"
$ws.Range("E15").Value = "  +0.77%  "

$ws.Range("D16").Value = "0.00001248"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("D17").Value = "1.600.22"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("D18").Value = "94.12"
$ws.Range("E18").Value = "  +2.44%  "

$ws.Range("D19").Value = "0.06924"
$ws.Range("E19").Value = "  +1.67%  "

$ws.Range("D20").Value = "18.10"
$ws.Range("E20").Value = "  -0.92%  "

$ws.Range("D21").Value = "6.530"
$ws.Range("E21").Value = "  -0.57%  "

$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.33%  "

$ws.Range("D23").Value = "12.89"
$ws.Range("E23").Value = "  -1.86%  "

$ws.Range("D24").Value = "23.215.30"
$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("D25").Value = "2.455"
$ws.Range("E25").Value = "  +4.29%  "

$ws.Range("D26").Value = "3.047"
$ws.Range("E26").Value = "  +7.24%  "

$ws.Range("D27").Value = "21.18"
$ws.Range("E27").Value = "  +0.62%  "

$ws.Range("D28").Value = "149.35"

$ws.Range("D29").Value = "5.273"
$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("D30").Value = "135.61"
$ws.Range("E30").Value = "  +1.24%  "

$ws.Range("E31").Value = "  +6.93%  "

$ws.Range("D32").Value = "6.716"
$ws.Range("E32").Value = "  -2.07%  "

$ws.Range("D33").Value = "1.776.41"
$ws.Range("E33").Value = "  -0.63%  "

$ws.Range("D34").Value = "0.9622"
$ws.Range("E34").Value = "  -0.82%  "

$ws.Range("D35").Value = "0.07472"
$ws.Range("E35").Value = "  -1.71%  "

$ws.Range("D36").Value = "10.31"
$ws.Range("E36").Value = "  -0.75%  "

$ws.Range("D37").Value = "0.02738"
$ws.Range("E37").Value = "  +1.17%  "

$ws.Range("D38").Value = "0.2528"
$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "6.109"
$ws.Range("E39").Value = "  -2.98%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "0.08771"
$ws.Range("E40").Value = "  -1.24%  "

$ws.Range("D41").Value = "1.384"
$ws.Range("E41").Value = "  +1.11%  "

$ws.Range("D42").Value = "0.7075"
$ws.Range("E42").Value = "  +0.52%  "

$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("D44").Value = "15.62"
$ws.Range("E44").Value = "  +1.64%  "

$ws.Range("D45").Value = "0.6529"
$ws.Range("E45").Value = "  -1.57%  "

$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("D47").Value = "2.313"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("D48").Value = "4.010"
$ws.Range("E48").Value = "  +0.50%  "

$ws.Range("D49").Value = "132.13"
$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").Value = "0.07918"
$ws.Range("E50").Value = "  -0.21%  "

$ws.Range("D51").Value = "1.201"
$ws.Range("E51").Value = "  -1.14%  "
